# Refresh NATMI LR-pair metrics (Ctf1-Il6st) with newly computed TPM-based
# values. Updates the derived-expression / specificity / edge-weight columns
# (G:J, M:T) in data rows 2-7; rows 5-7 also pick up an extra ligand-expressing
# cell (E) bumping the detection rate (F) to 1, which cascades into G:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: refreshed TPM-derived NATMI metrics
$ws.Range("G2").Value = 0.2076976666666666
$ws.Range("H2").Value = 0.6230929999999999
$ws.Range("I2").Value = 0.1089421893552267
$ws.Range("J2").Value = 0.1089421893552267
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 4.118164325296999
$ws.Range("R2").Value = 37.063478927673
$ws.Range("S2").Value = 0.01675954737287173
$ws.Range("T2").Value = 0.01675954737287172

# Row 3: refreshed TPM-derived NATMI metrics
$ws.Range("G3").Value = 0.2076976666666666
$ws.Range("H3").Value = 0.6230929999999999
$ws.Range("I3").Value = 0.1089421893552267
$ws.Range("J3").Value = 0.1089421893552267
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("S3").Value = 0.07194818805213708
$ws.Range("T3").Value = 0.07194818805213705

# Row 4: refreshed TPM-derived NATMI metrics
$ws.Range("G4").Value = 0.2076976666666666
$ws.Range("H4").Value = 0.6230929999999999
$ws.Range("I4").Value = 0.1089421893552267
$ws.Range("J4").Value = 0.1089421893552267
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 4.972020094777221
$ws.Range("R4").Value = 44.74818085299499
$ws.Range("S4").Value = 0.0202344539302179
$ws.Range("T4").Value = 0.0202344539302179

# Row 5: refreshed TPM-derived NATMI metrics
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.698796666666667
$ws.Range("H5").Value = 5.09639
$ws.Range("I5").Value = 0.8910578106447734
$ws.Range("J5").Value = 0.8910578106447733
$ws.Range("M5").Value = 19.827687
$ws.Range("N5").Value = 59.483061
$ws.Range("O5").Value = 0.1538389073329896
$ws.Range("P5").Value = 0.1538389073329896
$ws.Range("Q5").Value = 33.68320858331
$ws.Range("R5").Value = 303.14887724979
$ws.Range("S5").Value = 0.1370793599601179
$ws.Range("T5").Value = 0.1370793599601179

# Row 6: refreshed TPM-derived NATMI metrics
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.698796666666667
$ws.Range("H6").Value = 5.09639
$ws.Range("I6").Value = 0.8910578106447734
$ws.Range("J6").Value = 0.8910578106447733
$ws.Range("O6").Value = 0.6604253914664442
$ws.Range("P6").Value = 0.6604253914664441
$ws.Range("Q6").Value = 144.6009114347633
$ws.Range("R6").Value = 1301.40820291287
$ws.Range("S6").Value = 0.5884772034143072
$ws.Range("T6").Value = 0.588477203414307

# Row 7: refreshed TPM-derived NATMI metrics
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.698796666666667
$ws.Range("H7").Value = 5.09639
$ws.Range("I7").Value = 0.8910578106447734
$ws.Range("J7").Value = 0.8910578106447733
$ws.Range("M7").Value = 23.93873833333333
$ws.Range("N7").Value = 71.816215
$ws.Range("O7").Value = 0.1857357012005663
$ws.Range("P7").Value = 0.1857357012005663
$ws.Range("Q7").Value = 40.66704888487222
$ws.Range("R7").Value = 366.00343996385
$ws.Range("S7").Value = 0.1655012472703484
$ws.Range("T7").Value = 0.1655012472703484

